# Fix the "Geladeira - PDV" ordering typo in the SOVI SSD GELADO KO STORE TYPES
# cell of the KPIs template, and clear the leftover grey "needs review" highlight
# now that the value has been corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPIs")

$oldText = "Balcão Refrigerado – PDV, Geladeira – CONC, Geladeira – KO, GELADEIRA – PDV"
$newText = "Balcão Refrigerado – PDV, Geladeira – CONC, Geladeira – PDV, Geladeira – KO"

$rows = @(5, 10, 13, 16)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value -eq $oldText) {
        $cell.Value = $newText
    }
    # Remove the grey "pending fix" fill now that the text is correct.
    $cell.Interior.Pattern = -4142
    $cell.Interior.ColorIndex = -4142
}

# Row 5 no longer needs the extra height the highlighted text required.
$ws.Rows.Item(5).RowHeight = 13.8

# Restore the normal working selection near the top of the sheet.
$ws.Activate()
$ws.Range("A8").Select()
